# "Generate Report for Handback" - update the localization status report:
#   - Overview sheet: status text for zh-cn/de-de moves from "Ready for handoff"
#     to "Handed back: in sync with en-US"
#   - zh-cn / de-de detail sheets: Status text updated the same way, the
#     "Latest Handback DateTime" is refreshed, and the stale "Error Detail"
#     message (handback file was out of date) is cleared now that the
#     handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-13 13:21:34"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-10-13 13:21:50"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(16).ColumnWidth = 12.8
